# Improved Tasks 1 and 2 Backend
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# New table contents (rows 2-11, columns B:D).
# Row 2 is the header (kept as-is). Rows 3-11 are the service rows.
# ---------------------------------------------------------------

$data = @(
    @("News",               "Responsável por oferecer as noticias aos utilizadores",               "Lista de noticia"),
    @("Stats",               "Responsável por oferecer as estatisticas aos utilizadores",           "Lista de estatisticas positivas e negativas"),
    @("Market",              "Responsável pelas funcionalidades de mercados/ mercados especificos.","Supermercados, Categorias de carnes"),
    @("Product",             "Responsável por todas as funcionalidades relacionadas com as carnes", "Carnes"),
    @("Feedback",            "Responsável pelo feedback dos utilizadores",                          "Feedback"),
    @("Product In Market",   "Responsável por relacionar os produtos com os supermercados",         "Supermercados"),
    @("User",                "Responsável pela autenticação",                                       "Utilizador"),
    @("App Config",          "Responsável pela configuração da aplicação",                          "Definições"),
    @("Saved",               "Responsável pelas funcionalidades de objetos guardados",              "Produtos, Utilizador")
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]

    # Column B: centered (new style)
    $ws.Cells.Item($r, 2).HorizontalAlignment = -4108  # xlCenter
    $ws.Cells.Item($r, 2).VerticalAlignment = -4108    # xlCenter

    # Columns C and D: wrap text
    $ws.Cells.Item($r, 3).WrapText = $true
    $ws.Cells.Item($r, 4).WrapText = $true
}

# Row 3 / Row 9 : column C keeps no wrap (short text fits on one line)
$ws.Cells.Item(3, 3).WrapText = $false
$ws.Cells.Item(3, 4).WrapText = $false
$ws.Cells.Item(4, 3).WrapText = $false
$ws.Cells.Item(9, 4).WrapText = $false

# ---------------------------------------------------------------
# New column A (narrow spacer column)
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 5.88671875

# ---------------------------------------------------------------
# Sheet view: zoom + selection
# ---------------------------------------------------------------
$excel.ActiveWindow.Zoom = 100
$ws.Range("C8").Select()
